# Cover-Letter-Syed_Taha_Rezwan.docx edit script
#
# This applies three textual edits, each of which splits previously
# single/merged runs into multiple runs with identical formatting
# (matching the shape the author's Word session produced):
#
#   1. Header name "Syed " / "Taha Rezwan"      -> "Syed Taha R" / "i" / "zwan"
#   2. Body: "... candidate for this internship. I am keen ..."
#                                                -> "... this " / "job" / ". I am keen ..."
#   3. Signature "Syed Taha Rezwan "             -> "Syed Taha R" / "i" / "zwan "
#
# Notes on this harness's Range semantics (discovered empirically):
#   - Range.Text = "..." on a range whose bounds fall inside a run splits
#     that run but drops an explicitly-empty <w:rPr/> down to "no rPr"
#     unless some Font property is subsequently toggled on the resulting
#     sub-range (this also forces the split to "stick" rather than being
#     silently re-coalesced with an identically-formatted neighbour).
#   - The reliable pattern is: set the text, then toggle a Font property
#     off/on (or on/off) back to a neutral value so the run's rPr element
#     materializes / the split is preserved, without changing the visible
#     formatting.
#   - A `Paragraph.Range` object is "elastic" (tracks live positions in a
#     way that breaks a subsequent same-range `.Text =`), so we always
#     re-derive a plain Range via $d.Range(start, end) before mutating.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: header name paragraph, "Syed " + "Taha Rezwan" (sz=40 both)
#         -> "Syed Taha R" + "i" + "zwan" (sz=40 all three)
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$pr1 = $p1.Range
$head = $d.Range($pr1.Start, $pr1.End)
$head.Text = "Syed Taha R"
# force-split / keep explicit rPr by round-tripping the font size
$head.Font.Size = 99
$head.Font.Size = 20

$hIStart = $head.End
$head.InsertAfter("i")
$hI = $d.Range($hIStart, $hIStart + 1)
$hI.Font.Size = 99
$hI.Font.Size = 20

$hZStart = $hI.End
$hI.InsertAfter("zwan")
$hZ = $d.Range($hZStart, $hZStart + 4)
$hZ.Font.Size = 99
$hZ.Font.Size = 20

# ---------------------------------------------------------------------
# Edit 2: body paragraph, "... a strong candidate for this internship.
#         I am keen ..." -> "... this " + "job" + ". I am keen ..."
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$pr6 = $p6.Range
$b6Start = $pr6.Start
$b6End = $pr6.End
$body = $d.Range($b6Start, $b6End)
$bodyText = $body.Text
$hit = $bodyText.IndexOf("internship. I am keen")
$absHit = $b6Start + $hit

# "internship" (10 chars) -> "job"
$mid = $d.Range($absHit, $absHit + 10)
$mid.Text = "job"
$mid.Font.Bold = 1
$mid.Font.Bold = 0

# tail run: from end of "job" through end of (now-shorter) paragraph
$p6again = $d.Paragraphs.Item(6)
$b6NewEnd = $p6again.Range.End
$tail = $d.Range($mid.End, $b6NewEnd)
$tail.Font.Bold = 1
$tail.Font.Bold = 0

# lead run: restore its explicit (empty) rPr, which the mid-run edit
# above stripped from the untouched leading text
$lead = $d.Range($b6Start, $absHit)
$lead.Font.Bold = 1
$lead.Font.Bold = 0

# ---------------------------------------------------------------------
# Edit 3: signature paragraph, "Syed Taha Rezwan "
#         -> "Syed Taha R" + "i" + "zwan "
# ---------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$pr9 = $p9.Range
$sig = $d.Range($pr9.Start, $pr9.End)
$sig.Text = "Syed Taha R"
$sig.Font.Bold = 1
$sig.Font.Bold = 0

$sIStart = $sig.End
$sig.InsertAfter("i")
$sI = $d.Range($sIStart, $sIStart + 1)
$sI.Font.Bold = 1
$sI.Font.Bold = 0

$sZStart = $sI.End
$sI.InsertAfter("zwan ")
$sZ = $d.Range($sZStart, $sZStart + 5)
$sZ.Font.Bold = 1
$sZ.Font.Bold = 0

Write-Host "Edits applied."
